$d = $word.ActiveDocument

function Retype-Paragraph($paraIndex, $newText) {
    # Deleting the whole paragraph range (including its end-of-paragraph
    # mark) and re-inserting the text + a fresh paragraph mark rebuilds
    # the paragraph from a single clean run, which drops any stray
    # leftover proofing markers (w:proofErr) that a plain Range/Find
    # text replacement would otherwise leave orphaned.
    $p = $d.Paragraphs.Item($paraIndex).Range
    $start = $p.Start
    $full = $d.Range($start, $p.End)
    $full.Delete()
    $d.Range($start, $start).InsertBefore($newText + "`r")
}

Retype-Paragraph 1 "Daniela castaño garcia"
Retype-Paragraph 2 "Esto es una prueba para entender elementos de git."

# Insert the new paragraph after the "practica" paragraph (paragraph 3)
# and before the trailing (space-only) paragraph.
$p3 = $d.Paragraphs.Item(3).Range
$d.Range($p3.End, $p3.End).InsertAfter("Este es una nueva modificación para poner en práctica.`r")
